$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New mutual information values for column B (rows 2..91), replacing the
# previous run's results ("run with 10000 mutual information values").
$values = @(
5, 5, 4, 1, 4, 5, 3, 4, 5, 2, 2, 4, 3, 5, 5, 5, 2, 4, 5, 4, 5, 2, 2, 2, 2, 1, 5, 1, 3, 5, 5, 4, 4, 4, 2, 2, 3, 2, 4, 5, 2, 2, 5, 5, 4, 4, 2, 1, 5, 3, 2, 3, 4, 2, 4, 2, 3, 5, 1, 4, 2, 1, 4, 2, 5, 5, 4, 1, 5, 3, 5, 2, 4, 4, 5, 2, 5, 5, 3, 5, 3, 4, 5, 5, 5, 4, 4, 4, 2, 3
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
